$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.800.87'
$ws.Range("E2").Value = '  +2.02%  '
$ws.Range("D3").Value = '3.087.59'
$ws.Range("E3").Value = '  +4.99%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.50%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.082.24'
$ws.Range("E8").Value = '  +4.91%  '
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.61'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.156'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.482'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.07%  '
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").Value = '3.599.91'
$ws.Range("E16").Value = '  +5.03%  '
$ws.Range("D17").Value = '66.803.05'
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.37%  '
$ws.Range("D19").Value = '3.088.06'
$ws.Range("E19").Value = '  +5.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("E21").Value = '  +4.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.713'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.58%  '
$ws.Range("E23").Value = '  +3.34%  '
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("E25").Value = '  +6.12%  '
$ws.Range("E26").Value = '  +7.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.76%  '
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.98%  '
$ws.Range("E34").Value = '  +3.49%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +3.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.89'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.18%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.318'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.00%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.28'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.39%  '
$ws.Range("E44").Value = '  -1.30%  '
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '382.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("D47").Value = '2.781.59'
$ws.Range("E47").Value = '  +2.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.28%  '
$ws.Range("E51").Value = '  +2.07%  '
